# Applies the 2026-01-02 "Jogos do Dia" update:
#  - tweak several odds in the existing rows 2, 4 and 5
#  - insert a new fixture row (Cypriot 1st Division) at row 6, which pushes
#    the former rows 6-9 down to rows 7-10
#  - refresh the odds for the (now shifted) rows 7-10 to their new values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($ws, $r, $c, $val) {
    $ws.Cells.Item($r, $c).Value = $val
}

# The "Date" column holds plain text like "2026-01-02"; writing that string
# straight into .Value lets Excel's smart-typing turn it into a real date
# serial number. Forcing a text format first keeps it as text, then the
# style is reset to Normal so no stray number-format sticks to the cell.
function Set-DateCell($ws, $r, $c, $val) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Set-DataRow($ws, $rowNum, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 1
        if ($col -eq 2) {
            Set-DateCell $ws $rowNum $col $values[$i]
        } else {
            Set-Cell $ws $rowNum $col $values[$i]
        }
    }
}

# --- Update existing rows 2, 4, 5 (values per diff) ---
Set-Cell $ws 2 6 1.87
Set-Cell $ws 2 9 4.5
Set-Cell $ws 2 21 2.3
Set-Cell $ws 2 22 1.28
Set-Cell $ws 2 30 17.5
Set-Cell $ws 2 31 50

Set-Cell $ws 4 6 2.82
Set-Cell $ws 4 7 3
Set-Cell $ws 4 10 1.09
Set-Cell $ws 4 11 470
Set-Cell $ws 4 14 1.34
Set-Cell $ws 4 16 1.34
Set-Cell $ws 4 23 1.5
Set-Cell $ws 4 24 20
Set-Cell $ws 4 25 18
Set-Cell $ws 4 26 26
Set-Cell $ws 4 27 95
Set-Cell $ws 4 28 24
Set-Cell $ws 4 29 10.5
Set-Cell $ws 4 30 22
Set-Cell $ws 4 31 90
Set-Cell $ws 4 32 48
Set-Cell $ws 4 33 34
Set-Cell $ws 4 34 50
Set-Cell $ws 4 35 160
Set-Cell $ws 4 36 260
Set-Cell $ws 4 37 200
Set-Cell $ws 4 38 290
Set-Cell $ws 4 39 600
Set-Cell $ws 4 40 460
Set-Cell $ws 4 41 150

Set-Cell $ws 5 13 1.05

# --- Insert new row at position 6 (shifts rows 6-9 down to 7-10) ---
$ws.Rows("6:6").Insert()

# --- Set full row content for rows 6 through 10 ---
$row6 = @('Cypriot 1st Division', '2026-01-02', '14:00:00', 'Omonia FC Aradippou', 'Digenis Ypsona', 2.36, 2.88, 3.15, 4, 3, 3.5, 1.4, 1.01, 2.4, 1.01, 1.55, 2.1, 1.18, 3.7, 1.04, 1.04, 1.33, 1.53, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000)
Set-DataRow $ws 6 $row6

$row7 = @('Saudi Professional League', '2026-01-02', '14:30:00', 'Al Ahli', 'Al Nassr', 4.7, 5.6, 1.57, 1.71, 4.6, 5.6, 1.01, 1.03, 2.68, 1.15, 2.68, 1.47, 1.58, 2.02, 1.5, 2.16, 2.4, 1.22, 980, 980, 18.5, 980, 38, 980, 980, 21, 65, 980, 980, 980, 1000, 75, 70, 100, 1000, 1000)
Set-DataRow $ws 7 $row7

$row8 = @('French Ligue 1', '2026-01-02', '16:45:00', 'Toulouse', 'Lens', 2.9, 2.96, 2.66, 2.72, 3.4, 3.5, 0, 1.08, 3.7, 1.34, 1.9, 2, 1.35, 3.6, 1.78, 2.2, 0, 0, 13.5, 11, 22, 42, 12, 7.6, 12.5, 32, 23, 13, 16.5, 42, 46, 38, 46, 100, 32, 26)
Set-DataRow $ws 8 $row8

$row9 = @('Italian Serie A', '2026-01-02', '16:45:00', 'Cagliari', 'AC Milan', 7, 7.2, 1.58, 1.59, 4.4, 4.5, 0, 1.07, 3.85, 1.32, 1.98, 1.99, 1.37, 3.5, 2.06, 1.87, 0, 0, 15.5, 7.8, 8.800000000000001, 14.5, 20, 9.6, 10.5, 17.5, 60, 26, 27, 42, 260, 1000, 120, 190, 180, 9.4)
Set-DataRow $ws 9 $row9

$row10 = @('Spanish La Liga', '2026-01-02', '17:00:00', 'Rayo Vallecano', 'Getafe', 2.24, 2.28, 4.4, 4.6, 3, 3.05, 0, 1.17, 2.3, 1.75, 1.4, 3.3, 1.14, 7.6, 2.64, 1.58, 0, 0, 7, 9.800000000000001, 32, 160, 6, 7.4, 22, 110, 10.5, 13.5, 34, 170, 30, 38, 95, 430, 40, 240)
Set-DataRow $ws 10 $row10

Write-Host "Edit complete"